$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 462.42856
$ws.Cells.Item(19, 9).Value = 535.625
$ws.Cells.Item(19, 10).Value = 364.83334
$ws.Cells.Item(19, 11).Value = 535.625
$ws.Cells.Item(19, 12).Value = 364.83334
$ws.Cells.Item(19, 13).Value = -360.625
$ws.Cells.Item(19, 14).Value = -714.83334
$ws.Cells.Item(28, 8).Value = 622
$ws.Cells.Item(28, 9).Value = 725.8
$ws.Cells.Item(28, 11).Value = 725.8
$ws.Cells.Item(28, 13).Value = -240.8
$ws.Cells.Item(55, 8).Value = 130.5
$ws.Cells.Item(55, 9).Value = 203.6
$ws.Cells.Item(55, 11).Value = 203.6
$ws.Cells.Item(55, 13).Value = 10.40000000000001
$ws.Cells.Item(62, 8).Value = 6465.5835
$ws.Cells.Item(62, 9).Value = 4598.8
$ws.Cells.Item(62, 11).Value = 4598.8
$ws.Cells.Item(62, 13).Value = -3974.8
$ws.Cells.Item(65, 8).Value = 6465.5835
$ws.Cells.Item(65, 9).Value = 4598.8
$ws.Cells.Item(65, 11).Value = 22994
$ws.Cells.Item(65, 13).Value = -19874
$ws.Cells.Item(70, 8).Value = 2698.75
$ws.Cells.Item(70, 9).Value = 2597.8572
$ws.Cells.Item(70, 10).Value = 2840
$ws.Cells.Item(70, 11).Value = 7793.571599999999
$ws.Cells.Item(70, 12).Value = 8520
$ws.Cells.Item(70, 13).Value = -7523.571599999999
$ws.Cells.Item(70, 14).Value = -9060
$ws.Cells.Item(73, 8).Value = 2698.75
$ws.Cells.Item(73, 9).Value = 2597.8572
$ws.Cells.Item(73, 10).Value = 2840
$ws.Cells.Item(73, 11).Value = 7793.571599999999
$ws.Cells.Item(73, 12).Value = 8520
$ws.Cells.Item(73, 13).Value = -6857.571599999999
$ws.Cells.Item(73, 14).Value = -10392
$ws.Cells.Item(74, 8).Value = 9625
$ws.Cells.Item(74, 10).Value = 9625
$ws.Cells.Item(74, 12).Value = 9625
$ws.Cells.Item(74, 14).Value = -11497
$ws.Cells.Item(77, 8).Value = 9625
$ws.Cells.Item(77, 10).Value = 9625
$ws.Cells.Item(77, 12).Value = 48125
$ws.Cells.Item(77, 14).Value = -57485
$ws.Cells.Item(80, 8).Value = 499.125
$ws.Cells.Item(80, 10).Value = 385.875
$ws.Cells.Item(80, 12).Value = 1157.625
$ws.Cells.Item(80, 14).Value = -3153.625
$ws.Cells.Item(83, 8).Value = 499.125
$ws.Cells.Item(83, 10).Value = 385.875
$ws.Cells.Item(83, 12).Value = 3472.875
$ws.Cells.Item(83, 14).Value = -13456.875
$ws.Cells.Item(136, 8).Value = 92316.22
$ws.Cells.Item(136, 10).Value = 92316.22
$ws.Cells.Item(136, 12).Value = 92316.22
$ws.Cells.Item(136, 14).Value = -102516.22
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1526.25
$ws.Cells.Item(45, 9).Value = 1352.2858
$ws.Cells.Item(45, 11).Value = 1352.2858
$ws.Cells.Item(45, 13).Value = -975.2858000000001
$ws.Cells.Item(74, 8).Value = 8072.579
$ws.Cells.Item(74, 9).Value = 2788.4375
$ws.Cells.Item(74, 10).Value = 36254.668
$ws.Cells.Item(74, 11).Value = 2788.4375
$ws.Cells.Item(74, 12).Value = 36254.668
$ws.Cells.Item(74, 13).Value = -1914.4375
$ws.Cells.Item(74, 14).Value = -38002.668
$ws.Cells.Item(77, 8).Value = 8072.579
$ws.Cells.Item(77, 9).Value = 2788.4375
$ws.Cells.Item(77, 10).Value = 36254.668
$ws.Cells.Item(77, 11).Value = 13942.1875
$ws.Cells.Item(77, 12).Value = 181273.34
$ws.Cells.Item(77, 13).Value = -9574.1875
$ws.Cells.Item(77, 14).Value = -190009.34
$ws.Cells.Item(132, 8).Value = 773180.9399999999
$ws.Cells.Item(132, 9).Value = 5120.7188
$ws.Cells.Item(132, 11).Value = 15362.1564
$ws.Cells.Item(132, 13).Value = -12832.1564
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 14295.375
$ws.Cells.Item(22, 9).Value = 14295.375
$ws.Cells.Item(22, 11).Value = 14295.375
$ws.Cells.Item(22, 13).Value = -14122.375
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 13).ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(9, 8).Value = 99000
$ws.Cells.Item(9, 10).Value = 99000
$ws.Cells.Item(9, 12).Value = 99000
$ws.Cells.Item(9, 14).Value = -99336
$ws.Cells.Item(31, 8).Value = 81972.13
$ws.Cells.Item(31, 9).Value = 106548
$ws.Cells.Item(31, 11).Value = 106548
$ws.Cells.Item(31, 13).Value = -106253
$ws.Cells.Item(34, 8).Value = 81972.13
$ws.Cells.Item(34, 9).Value = 106548
$ws.Cells.Item(34, 11).Value = 106548
$ws.Cells.Item(34, 13).Value = -106346
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 14).ClearContents()
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 14).ClearContents()
$ws.Cells.Item(88, 8).Value = 32500
$ws.Cells.Item(88, 10).Value = 32500
$ws.Cells.Item(88, 12).Value = 32500
$ws.Cells.Item(88, 14).Value = -33312
$ws.Cells.Item(91, 8).Value = 32500
$ws.Cells.Item(91, 10).Value = 32500
$ws.Cells.Item(91, 12).Value = 32500
$ws.Cells.Item(91, 14).Value = -35308
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 1833.1666
$ws.Cells.Item(8, 9).Value = 1833.1666
$ws.Cells.Item(8, 11).Value = 5499.4998
$ws.Cells.Item(8, 13).Value = -5360.4998
$ws.Cells.Item(23, 8).Value = 724
$ws.Cells.Item(23, 9).Value = 169.25
$ws.Cells.Item(23, 11).Value = 507.75
$ws.Cells.Item(23, 13).Value = -272.75
$ws.Cells.Item(26, 8).Value = 1832.5
$ws.Cells.Item(26, 9).Value = 2203.8
$ws.Cells.Item(26, 10).Value = 1567.2858
$ws.Cells.Item(26, 11).Value = 6611.400000000001
$ws.Cells.Item(26, 12).Value = 4701.857400000001
$ws.Cells.Item(26, 13).Value = -6323.400000000001
$ws.Cells.Item(26, 14).Value = -5277.857400000001
$ws.Cells.Item(39, 8).Value = 6540.6665
$ws.Cells.Item(39, 9).Value = 3722
$ws.Cells.Item(39, 10).Value = 7950
$ws.Cells.Item(39, 11).Value = 11166
$ws.Cells.Item(39, 12).Value = 23850
$ws.Cells.Item(39, 13).Value = -10872
$ws.Cells.Item(39, 14).Value = -24438
$ws.Cells.Item(62, 8).Value = 3156.8572
$ws.Cells.Item(62, 9).Value = 3156.8572
$ws.Cells.Item(62, 11).Value = 9470.571599999999
$ws.Cells.Item(62, 13).Value = -8784.571599999999
$ws.Cells.Item(65, 8).Value = 3156.8572
$ws.Cells.Item(65, 9).Value = 3156.8572
$ws.Cells.Item(65, 11).Value = 28411.7148
$ws.Cells.Item(65, 13).Value = -24979.7148
$ws.Cells.Item(107, 8).Value = 5231.9165
$ws.Cells.Item(107, 10).Value = 7559.25
$ws.Cells.Item(107, 12).Value = 22677.75
$ws.Cells.Item(107, 14).Value = -26517.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 4816.476
$ws.Cells.Item(113, 10).Value = 4414.1665
$ws.Cells.Item(113, 12).Value = 4414.1665
$ws.Cells.Item(113, 14).Value = -8754.166499999999
$ws.Cells.Item(126, 8).Value = 7269
$ws.Cells.Item(126, 9).Value = 8711.200000000001
$ws.Cells.Item(126, 11).Value = 26133.6
$ws.Cells.Item(126, 13).Value = -23663.6
$ws.Cells.Item(132, 8).Value = 1070759.5
$ws.Cells.Item(132, 9).Value = 6665.222
$ws.Cells.Item(132, 10).Value = 2438880.8
$ws.Cells.Item(132, 11).Value = 19995.666
$ws.Cells.Item(132, 12).Value = 7316642.399999999
$ws.Cells.Item(132, 13).Value = -17465.666
$ws.Cells.Item(132, 14).Value = -7321702.399999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 499.5
$ws.Cells.Item(22, 9).Value = 499.5
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 499.5
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -204.5
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(27, 8).Value = 499.5
$ws.Cells.Item(27, 9).Value = 499.5
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 11).Value = 499.5
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 13).Value = -392.5
$ws.Cells.Item(27, 14).ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(58, 8).Value = 9085
$ws.Cells.Item(58, 9).Value = 9085
$ws.Cells.Item(58, 11).Value = 9085
$ws.Cells.Item(58, 13).Value = -8777
$ws.Cells.Item(70, 8).Value = 28997.5
$ws.Cells.Item(70, 10).Value = 28997.5
$ws.Cells.Item(70, 12).Value = 28997.5
$ws.Cells.Item(70, 14).Value = -29627.5
$ws.Cells.Item(73, 8).Value = 28997.5
$ws.Cells.Item(73, 10).Value = 28997.5
$ws.Cells.Item(73, 12).Value = 28997.5
$ws.Cells.Item(73, 14).Value = -31181.5
